$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two weekly data blocks (rows 2-3 <-> rows 4-5) for columns
# D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), S (Precio $/Kg).
# Column L (Calidad: Primera/Segunda) stays fixed per row.

$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $topRange = "$col" + "2:" + "$col" + "3"
    $bottomRange = "$col" + "4:" + "$col" + "5"

    $topValues = $ws.Range($topRange).Value2
    $bottomValues = $ws.Range($bottomRange).Value2

    $ws.Range($topRange).Value2 = $bottomValues
    $ws.Range($bottomRange).Value2 = $topValues
}
